# Updates cryptos list price/volume figures (and restores the Mantle /
# EnergySwap row ordering) to match the latest scrape.
#
# Columns D (Price) and E (Volume(1h)) are stored as plain text in this
# sheet, even when the text looks like a number (e.g. "215.30"). Assigning
# such a string straight to .Value would make Excel auto-detect it as a
# number, so for any new Price value that parses as a plain number we
# briefly mark the cell as Text ("@"), write the value, then restore the
# cell's original ("Normal") style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "25.803.49"
$ws.Cells.Item(2, 5).Value = "  -0.13%  "
$ws.Cells.Item(3, 4).Value = "1.636.00"
$ws.Cells.Item(3, 5).Value = "  -0.12%  "
$ws.Cells.Item(4, 5).Value = "  -0.09%  "
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "215.30"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.39%  "
$ws.Cells.Item(6, 5).Value = "  -0.81%  "
$ws.Cells.Item(7, 5).Value = "  -0.13%  "
$ws.Cells.Item(8, 5).Value = "  -0.16%  "
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0642"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.36%  "
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "19.87"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +1.17%  "
$ws.Cells.Item(11, 5).Value = "  +0.61%  "
$ws.Cells.Item(12, 5).Value = "  -0.88%  "
$ws.Cells.Item(13, 4).Value = "1.641.08"
$ws.Cells.Item(13, 5).Value = "  +0.20%  "
$ws.Cells.Item(14, 4).Value = "1.861.50"
$ws.Cells.Item(14, 5).Value = "  -0.18%  "
$ws.Cells.Item(15, 5).Value = "  -1.10%  "
$ws.Cells.Item(16, 4).Value = "0.0₃0775"
$ws.Cells.Item(16, 5).Value = "  +1.84%  "
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "63.06"
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.32%  "
$ws.Cells.Item(18, 4).Value = "25.812.86"
$ws.Cells.Item(18, 5).Value = "  -0.24%  "
$ws.Cells.Item(19, 5).Value = "  -0.15%  "
$ws.Cells.Item(20, 5).Value = "  +2.65%  "
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "193.89"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.98%  "
$ws.Cells.Item(22, 5).Value = "  +0.57%  "
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.17"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.91%  "
$ws.Cells.Item(24, 5).Value = "  -0.09%  "
$ws.Cells.Item(25, 5).Value = "  -0.69%  "
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "139.18"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.53%  "
$ws.Cells.Item(27, 5).Value = "  -4.75%  "
$ws.Cells.Item(28, 5).Value = "  +0.59%  "
$ws.Cells.Item(29, 5).Value = "  +0.22%  "
$ws.Cells.Item(30, 5).Value = "  +0.08%  "
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0497"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +1.65%  "
$ws.Cells.Item(33, 5).Value = "  +0.84%  "
$ws.Cells.Item(34, 5).Value = "  +2.51%  "
$ws.Cells.Item(35, 5).Value = "  +0.61%  "
$ws.Cells.Item(36, 5).Value = "  -0.73%  "
$ws.Cells.Item(37, 5).Value = "  -0.23%  "
$ws.Cells.Item(38, 5).Value = "  +0.23%  "
$ws.Cells.Item(39, 4).Value = "1.108.32"
$ws.Cells.Item(39, 5).Value = "  -1.91%  "
$ws.Cells.Item(40, 5).Value = "  +0.25%  "
$ws.Cells.Item(41, 5).Value = "  +0.72%  "
$ws.Cells.Item(42, 5).Value = "  +0.95%  "
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "99.24"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +1.55%  "
$ws.Cells.Item(44, 5).Value = "  +0.03%  "
$ws.Cells.Item(45, 4).Value = "0.0₆0111"
$ws.Cells.Item(45, 5).Value = "  -1.64%  "
$ws.Cells.Item(46, 5).Value = "  +0.22%  "
$ws.Cells.Item(47, 5).Value = "  +12.93%  "
$ws.Cells.Item(48, 2).Value = "Mantle"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.418"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -6.04%  "
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.70"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -0.24%  "
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0503"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.36%  "
$ws.Cells.Item(51, 5).Value = "  +0.03%  "
